# B6-PowerPoint.pptx edit
#
# 1) Three tables (on the slides that hold the "Tasks" comparison tables)
#    get their table style switched from the locally-defined default style
#    {2A8C8D4F-A319-4E8F-9F79-D107947C3321} to the built-in style
#    {F597ACA8-54AF-4C73-9E65-B428E045F97C}.
#
# 2) The deck's theme colour scheme changes from the "Integral / Red Violet"
#    palette to the plain "Office" palette (the font scheme and format
#    scheme are identical between the two themes - only the 12 theme
#    colours differ), so the presentation's overall design switches from
#    the red/violet accent palette to the standard Office blue/orange one.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the three tables
# ---------------------------------------------------------------------
$newTableStyleId = "{F597ACA8-54AF-4C73-9E65-B428E045F97C}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId, $false)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the theme colours: Integral/Red Violet -> Office
# ---------------------------------------------------------------------
function ConvertHexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorScheme.Colors(1..12):
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertHexToRgb($officeThemeColors[$i - 1])
}
